# Updated cryptos list with GitHub Actions
#
# Refreshes the scraped price / 1h volume-change figures on the
# "cryptos" worksheet and fixes the ordering of two rows whose
# underlying coins had been swapped (FraxShare <-> ApeXProtocol).
#
# Many of the "Price" values look like plain decimal numbers
# (e.g. "96.74"), but the source data treats the whole column as
# text (some rows use thousands separators, e.g. "42.761.90",
# which are not valid numbers). To avoid Excel silently converting
# the numeric-looking ones into floating point numbers (which can
# introduce binary rounding artifacts such as 96.7400000000001),
# we force the cell to Text format before assigning the value, and
# restore the cell's normal style afterwards so formatting stays
# untouched.
function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
Set-TextValue $ws "D2" "42.761.90"
$ws.Range("E2").Value = "  -0.41%  "

# --- Row 3: Ethereum ---
Set-TextValue $ws "D3" "2.294.55"
$ws.Range("E3").Value = "  -0.05%  "

# --- Row 4: TetherUSD ---
Set-TextValue $ws "D4" "0.999"
$ws.Range("E4").Value = "  -0.06%  "

# --- Row 5: BNB ---
Set-TextValue $ws "D5" "303.81"
$ws.Range("E5").Value = "  +1.37%  "

# --- Row 6: Solana ---
Set-TextValue $ws "D6" "96.74"
$ws.Range("E6").Value = "  -0.70%  "

# --- Row 7: XRP ---
Set-TextValue $ws "D7" "0.505"
$ws.Range("E7").Value = "  -1.97%  "

# --- Row 8: USDC ---
$ws.Range("E8").Value = "  +0.01%  "

# --- Row 10: Avalanche ---
$ws.Range("E10").Value = "  -1.82%  "

# --- Row 11: Dogecoin ---
$ws.Range("E11").Value = "  -0.71%  "

# --- Row 12: Chainlink ---
Set-TextValue $ws "D12" "18.79"
$ws.Range("E12").Value = "  +6.18%  "

# --- Row 13: TRON ---
Set-TextValue $ws "D13" "0.120"
$ws.Range("E13").Value = "  +2.12%  "

# --- Row 14: Polkadot ---
$ws.Range("E14").Value = "  +1.05%  "

# --- Row 15: WrappedliquidstakedEther2.0 ---
Set-TextValue $ws "D15" "2.651.89"
$ws.Range("E15").Value = "  -0.11%  "

# --- Row 16: WrappedEther ---
Set-TextValue $ws "D16" "2.291.16"
$ws.Range("E16").Value = "  -1.09%  "

# --- Row 17: Polygon ---
$ws.Range("E17").Value = "  -0.28%  "

# --- Row 18: WrappedBTC ---
Set-TextValue $ws "D18" "42.687.46"
$ws.Range("E18").Value = "  -0.44%  "

# --- Row 19: InternetComputer(DFINITY) ---
Set-TextValue $ws "D19" "12.83"
$ws.Range("E19").Value = "  +1.71%  "

# --- Row 20: ShibaInu ---
$ws.Range("E20").Value = "  -1.31%  "

# --- Row 21: Uniswap ---
Set-TextValue $ws "D21" "6.00"
$ws.Range("E21").Value = "  -1.64%  "

# --- Row 22: Litecoin ---
Set-TextValue $ws "D22" "67.23"
$ws.Range("E22").Value = "  -1.18%  "

# --- Row 23: BitcoinCash ---
Set-TextValue $ws "D23" "236.02"
$ws.Range("E23").Value = "  -2.26%  "

# --- Row 24: ImmutableX ---
$ws.Range("E24").Value = "  +0.67%  "

# --- Row 25: Dai ---
$ws.Range("E25").Value = "  +0.05%  "

# --- Row 26: PancakeSwap ---
Set-TextValue $ws "D26" "2.40"
$ws.Range("E26").Value = "  -1.28%  "

# --- Row 27: EthereumClassic ---
Set-TextValue $ws "D27" "25.04"
$ws.Range("E27").Value = "  -0.19%  "

# --- Row 28: Monero ---
Set-TextValue $ws "D28" "167.33"
$ws.Range("E28").Value = "  +0.47%  "

# --- Row 29: Toncoin ---
$ws.Range("E29").Value = "  +0.96%  "

# --- Row 30: Cosmos ---
$ws.Range("E30").Value = "  -0.46%  "

# --- Row 31: InjectiveProtocol ---
Set-TextValue $ws "D31" "33.00"
$ws.Range("E31").Value = "  +0.39%  "

# --- Row 32: FirstDigitalUSD ---
$ws.Range("E32").Value = "  +0.01%  "

# --- Row 33: Celestia ---
Set-TextValue $ws "D33" "17.88"
$ws.Range("E33").Value = "  +1.71%  "

# --- Row 34: Filecoin ---
$ws.Range("E34").Value = "  -0.66%  "

# --- Row 35: RenderToken ---
Set-TextValue $ws "D35" "4.47"
$ws.Range("E35").Value = "  -4.67%  "

# --- Row 36: WEMIXToken ---
$ws.Range("E36").Value = "  -1.68%  "

# --- Row 37: Hedera ---
$ws.Range("E37").Value = "  -0.39%  "

# --- Row 38: Kaspa ---
$ws.Range("E38").Value = "  -0.09%  "

# --- Row 39: ARBITRUM ---
$ws.Range("E39").Value = "  -0.62%  "

# --- Row 40: Stellar ---
$ws.Range("E40").Value = "  -1.02%  "

# --- Row 41: LidoDAOToken ---
$ws.Range("E41").Value = "  -2.38%  "

# --- Row 42: Maker ---
Set-TextValue $ws "D42" "1.993.73"
$ws.Range("E42").Value = "  -0.52%  "

# --- Row 43: VeChain ---
$ws.Range("E43").Value = "  -2.22%  "

# --- Row 44: EnergySwap ---
Set-TextValue $ws "D44" "18.44"
$ws.Range("E44").Value = "  +6.62%  "

# --- Rows 45-46: coins swapped places (FraxShare <-> ApeXProtocol) ---
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws "D45" "2.17"
$ws.Range("E45").Value = "  +0.68%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws "D46" "10.16"
$ws.Range("E46").Value = "  +0.16%  "

# --- Row 47: NEARProtocol ---
$ws.Range("E47").Value = "  +0.27%  "

# --- Row 48: HuobiToken ---
Set-TextValue $ws "D48" "2.90"
$ws.Range("E48").Value = "  -0.60%  "

# --- Row 49: MultiversX ---
Set-TextValue $ws "D49" "53.69"
$ws.Range("E49").Value = "  +0.55%  "

# --- Row 50: RocketPoolETH ---
Set-TextValue $ws "D50" "2.518.61"
$ws.Range("E50").Value = "  -0.21%  "

# --- Row 51: TrustWalletToken ---
$ws.Range("E51").Value = "  +1.58%  "

Write-Host "Applied cryptos list update"
